# Auto-generated script applying the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# "Balmung_Profits" price-refresh edit described by the commit diff.
# Each leve row's market-price-derived columns (H..N) are updated in place
# to reflect the new scraped prices; a few rows gain/lose a trailing cell.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 5561503.5
$ws.Range("I9").Value = 16667357
$ws.Range("K9").Value = 16667357
$ws.Range("M9").Value = -16667188
$ws.Range("H12").Value = 45454544
$ws.Range("I12").Value = 45454544
$ws.Range("K12").Value = 45454544
$ws.Range("M12").Value = -45454374
$ws.Range("H17").Value = 138566.52
$ws.Range("J17").Value = 141493.47
$ws.Range("L17").Value = 424480.41
$ws.Range("N17").Value = -424816.41
$ws.Range("H28").Value = 767.8889
$ws.Range("J28").Value = 999.5
$ws.Range("L28").Value = 999.5
$ws.Range("N28").Value = -1969.5
$ws.Range("H42").Value = 1240.3334
$ws.Range("J42").Value = 1971.4286
$ws.Range("L42").Value = 5914.2858
$ws.Range("N42").Value = -6374.2858
$ws.Range("H51").Value = 26371886
$ws.Range("J51").Value = 50002188
$ws.Range("L51").Value = 50002188
$ws.Range("N51").Value = -50003156
$ws.Range("H62").Value = 1963.25
$ws.Range("I62").Value = 1817.375
$ws.Range("J62").Value = 2109.125
$ws.Range("K62").Value = 1817.375
$ws.Range("L62").Value = 2109.125
$ws.Range("M62").Value = -1193.375
$ws.Range("N62").Value = -3357.125
$ws.Range("H65").Value = 1963.25
$ws.Range("I65").Value = 1817.375
$ws.Range("J65").Value = 2109.125
$ws.Range("K65").Value = 9086.875
$ws.Range("L65").Value = 10545.625
$ws.Range("M65").Value = -5966.875
$ws.Range("N65").Value = -16785.625
$ws.Range("H111").Value = 1125.7273
$ws.Range("I111").Value = 1121.1428
$ws.Range("J111").Value = 1133.75
$ws.Range("K111").Value = 3363.4284
$ws.Range("L111").Value = 3401.25
$ws.Range("M111").Value = -296.4284000000002
$ws.Range("N111").Value = -9535.25
$ws.Range("H116").Value = 12215.538
$ws.Range("I116").Value = 14597.2
$ws.Range("K116").Value = 14597.2
$ws.Range("M116").Value = -11155.2
$ws.Range("H118").Value = 689.0909
$ws.Range("I118").Value = 541.25
$ws.Range("K118").Value = 1623.75
$ws.Range("M118").Value = 33.25
$ws.Range("H132").Value = 32293.031
$ws.Range("I132").Value = 34312.6
$ws.Range("K132").Value = 102937.8
$ws.Range("M132").Value = -100407.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3062.7144
$ws.Range("I2").Value = 1073.1666
$ws.Range("K2").Value = 1073.1666
$ws.Range("M2").Value = -960.1666
$ws.Range("H37").Value = 72498.5
$ws.Range("J37").Value = 99999
$ws.Range("L37").Value = 99999
$ws.Range("N37").Value = -100545
$ws.Range("H45").Value = 54018.906
$ws.Range("I45").Value = 88283.086
$ws.Range("K45").Value = 88283.086
$ws.Range("M45").Value = -87906.086
$ws.Range("H61").Value = 1417554.5
$ws.Range("I61").Value = 33516.484
$ws.Range("K61").Value = 33516.484
$ws.Range("M61").Value = -33304.484
$ws.Range("H74").Value = 311521
$ws.Range("I74").Value = 1492.7046
$ws.Range("K74").Value = 1492.7046
$ws.Range("M74").Value = -618.7046
$ws.Range("H77").Value = 311521
$ws.Range("I77").Value = 1492.7046
$ws.Range("K77").Value = 7463.523
$ws.Range("M77").Value = -3095.523
$ws.Range("H110").Value = 1982.6666
$ws.Range("I110").Value = 1982.6666
$ws.Range("K110").Value = 1982.6666
$ws.Range("M110").Value = 62.33339999999998
$ws.Range("H114").Value = 80000
$ws.Range("J114").Value = 80000
$ws.Range("L114").Value = 80000
$ws.Range("N114").Value = -88678
$ws.Range("H116").Value = 3062.7144
$ws.Range("I116").Value = 1073.1666
$ws.Range("K116").Value = 1073.1666
$ws.Range("M116").Value = 1220.8334
$ws.Range("H136").Value = 1417554.5
$ws.Range("I136").Value = 33516.484
$ws.Range("K136").Value = 100549.452
$ws.Range("M136").Value = -97999.45199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3062.7144
$ws.Range("I3").Value = 1073.1666
$ws.Range("K3").Value = 1073.1666
$ws.Range("M3").Value = -959.1666
$ws.Range("H62").Value = 24000
$ws.Range("J62").Value = 24000
$ws.Range("L62").Value = 24000
$ws.Range("N62").Value = -25372
$ws.Range("H65").Value = 24000
$ws.Range("J65").Value = 24000
$ws.Range("L65").Value = 72000
$ws.Range("N65").Value = -78864
$ws.Range("H94").Value = 1641.8948
$ws.Range("I94").Value = 1103.4482
$ws.Range("K94").Value = 1103.4482
$ws.Range("M94").Value = -652.4482
$ws.Range("H107").Value = 6769.439
$ws.Range("I107").Value = 7978.2188
$ws.Range("J107").Value = 2471.5557
$ws.Range("K107").Value = 7978.2188
$ws.Range("L107").Value = 2471.5557
$ws.Range("M107").Value = -6058.2188
$ws.Range("N107").Value = -6311.5557
$ws.Range("H134").Value = 20456824
$ws.Range("I134").Value = 2029.2258
$ws.Range("K134").Value = 6087.6774
$ws.Range("M134").Value = -3552.6774

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H105").Value = 1794.375
$ws.Range("I105").Value = 1362.4615
$ws.Range("K105").Value = 1362.4615
$ws.Range("M105").Value = 384.5385000000001
$ws.Range("H132").Value = 25021.627
$ws.Range("I132").Value = 30956.176
$ws.Range("J132").Value = 2602.2222
$ws.Range("K132").Value = 92868.52799999999
$ws.Range("L132").Value = 7806.6666
$ws.Range("M132").Value = -90338.52799999999
$ws.Range("N132").Value = -12866.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 1476
$ws.Range("I57").Value = 1476
$ws.Range("K57").Value = 4428
$ws.Range("M57").Value = -3869
$ws.Range("H113").Value = 392.1111
$ws.Range("I113").Value = 384.5
$ws.Range("J113").Value = 393.43478
$ws.Range("K113").Value = 1153.5
$ws.Range("L113").Value = 1180.30434
$ws.Range("M113").Value = 1016.5
$ws.Range("N113").Value = -5520.30434
$ws.Range("H131").Value = 6496538.5
$ws.Range("I131").Value = 15153890
$ws.Range("J131").Value = 3525
$ws.Range("K131").Value = 45461670
$ws.Range("L131").Value = 10575
$ws.Range("M131").Value = -45456630
$ws.Range("N131").Value = -20655

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 76.46154
$ws.Range("J2").Value = 68.75
$ws.Range("L2").Value = 68.75
$ws.Range("N2").Value = -294.75
$ws.Range("H80").Value = 12766.625
$ws.Range("I80").Value = 11733.077
$ws.Range("J80").Value = 13988.091
$ws.Range("K80").Value = 11733.077
$ws.Range("L80").Value = 13988.091
$ws.Range("M80").Value = -10735.077
$ws.Range("N80").Value = -15984.091
$ws.Range("H83").Value = 12766.625
$ws.Range("I83").Value = 11733.077
$ws.Range("J83").Value = 13988.091
$ws.Range("K83").Value = 58665.38499999999
$ws.Range("L83").Value = 69940.455
$ws.Range("M83").Value = -53673.38499999999
$ws.Range("N83").Value = -79924.455
$ws.Range("H102").Value = 13159546
$ws.Range("I102").Value = 17242958
$ws.Range("J102").Value = 1889.4445
$ws.Range("K102").Value = 17242958
$ws.Range("L102").Value = 1889.4445
$ws.Range("M102").Value = -17241336
$ws.Range("N102").Value = -5133.4445
$ws.Range("H133").Value = 136461.67
$ws.Range("J133").Value = 136461.67
$ws.Range("L133").Value = 136461.67
$ws.Range("N133").Value = -146581.67

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6252.5713
$ws.Range("I7").Value = 2052.1333
$ws.Range("K7").Value = 2052.1333
$ws.Range("M7").Value = -1940.1333
$ws.Range("H16").Value = 1344.1818
$ws.Range("J16").Value = 3037
$ws.Range("L16").Value = 3037
$ws.Range("N16").Value = -3377
$ws.Range("H40").Value = 1564.1666
$ws.Range("I40").Value = 846.75
$ws.Range("K40").Value = 846.75
$ws.Range("M40").Value = -710.75
$ws.Range("H100").Value = 3750.1765
$ws.Range("I100").Value = 3783.7334
$ws.Range("J100").Value = 3498.5
$ws.Range("K100").Value = 3783.7334
$ws.Range("L100").Value = 3498.5
$ws.Range("M100").Value = -3242.7334
$ws.Range("N100").Value = -4580.5
$ws.Range("H126").Value = 6252.5713
$ws.Range("I126").Value = 2052.1333
$ws.Range("K126").Value = 6156.3999
$ws.Range("M126").Value = -3686.3999
$ws.Range("H132").Value = 3743.4443
$ws.Range("I132").Value = 3450.4
$ws.Range("K132").Value = 10351.2
$ws.Range("M132").Value = -7821.200000000001
$ws.Range("H136").Value = 1720.5769
$ws.Range("I136").Value = 2740.5293
$ws.Range("J136").Value = 1225.1714
$ws.Range("K136").Value = 8221.5879
$ws.Range("L136").Value = 3675.5142
$ws.Range("M136").Value = -5671.5879
$ws.Range("N136").Value = -8775.5142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3010.2122
$ws.Range("I132").Value = 3657.1
$ws.Range("J132").Value = 2728.9565
$ws.Range("K132").Value = 10971.3
$ws.Range("L132").Value = 8186.869499999999
$ws.Range("M132").Value = -8441.299999999999
$ws.Range("N132").Value = -13246.8695
$ws.Range("H136").Value = 43133.125
$ws.Range("I136").Value = 67274.8
$ws.Range("J136").Value = 2897
$ws.Range("K136").Value = 201824.4
$ws.Range("L136").Value = 8691
$ws.Range("M136").Value = -199274.4
$ws.Range("N136").Value = -13791
